$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-free inline edits.
# Columns A (CEDULA) and C (TELEFONO) carry a "Text" number format (style index 1)
# on data rows. When a numeric-looking string is re-entered into an already
# text-formatted cell it is kept as text (leading zeros preserved). Several
# phone numbers in the new data do NOT have a leading zero, and in the source
# workbook those ended up stored as real numbers (while still keeping the
# Text-style formatting index). To reproduce that exactly we reset the style
# to Normal first, write the numeric value, then restore the text number
# format - that sequence keeps the value numeric instead of re-coercing it.

# Row 2
$ws.Cells.Item(2,1).Value = "0961548488"
$ws.Cells.Item(2,2).Value = "ElenaParedes"
$ws.Cells.Item(2,3).Value = "0977845888"
$ws.Cells.Item(2,4).Value = 18
$ws.Cells.Item(2,5).Value = "SVPL543000,Gmys432000"

# Row 3
$ws.Cells.Item(3,1).Value = "0954872314"
$ws.Cells.Item(3,2).Value = "Carlos Jiménez"
$ws.Cells.Item(3,3).Value = "0987541123"
$ws.Cells.Item(3,4).Value = 25
$ws.Cells.Item(3,5).Value = "ADPR321000"

# Row 4 (C4 becomes a literal number 996312457, keeping the text style index)
$ws.Cells.Item(4,1).Value = "0912457896"
$ws.Cells.Item(4,2).Value = "María Torres"
$ws.Cells.Item(4,3).Style = "Normal"
$ws.Cells.Item(4,3).Value = 996312457
$ws.Cells.Item(4,3).NumberFormat = "@"
$ws.Cells.Item(4,4).Value = 32
$ws.Cells.Item(4,5).Value = "SVPL543000"

# Row 5
$ws.Cells.Item(5,1).Value = "0923654789"
$ws.Cells.Item(5,2).Value = "Juan Cedeño"
$ws.Cells.Item(5,3).Value = "0987123654"
$ws.Cells.Item(5,4).Value = 45
$ws.Cells.Item(5,5).Value = "HOGR554200"

# Row 6
$ws.Cells.Item(6,1).Value = "0932145789"
$ws.Cells.Item(6,2).Value = "Andrea Villalba"
$ws.Cells.Item(6,3).Value = "0978456321"
$ws.Cells.Item(6,4).Value = 29
$ws.Cells.Item(6,5).Value = "SGRC765000"

# Row 7
$ws.Cells.Item(7,1).Value = "0945873214"
$ws.Cells.Item(7,2).Value = "Pedro Suárez"
$ws.Cells.Item(7,3).Value = "0989658745"
$ws.Cells.Item(7,4).Value = 51
$ws.Cells.Item(7,5).Value = "SVPL543000"

# Row 8
$ws.Cells.Item(8,1).Value = "0956321478"
$ws.Cells.Item(8,2).Value = "Sofía Carrillo"
$ws.Cells.Item(8,3).Value = "0997412586"
$ws.Cells.Item(8,4).Value = 23
$ws.Cells.Item(8,5).Value = "SVPL543000"

# Row 9 (C9 becomes a literal number 984561237, keeping the text style index)
$ws.Cells.Item(9,1).Value = "0968745213"
$ws.Cells.Item(9,2).Value = "Luis Herrera"
$ws.Cells.Item(9,3).Style = "Normal"
$ws.Cells.Item(9,3).Value = 984561237
$ws.Cells.Item(9,3).NumberFormat = "@"
$ws.Cells.Item(9,4).Value = 38
$ws.Cells.Item(9,5).Value = "SGRC765000,ADPR321000"

# Row 10 (C10 becomes a literal number 998564123, keeping the text style index)
$ws.Cells.Item(10,1).Value = "0974512368"
$ws.Cells.Item(10,2).Value = "Daniela Morán"
$ws.Cells.Item(10,3).Style = "Normal"
$ws.Cells.Item(10,3).Value = 998564123
$ws.Cells.Item(10,3).NumberFormat = "@"
$ws.Cells.Item(10,4).Value = 30
$ws.Cells.Item(10,5).Value = "SRVJ987000"

# Row 11 (new row). A11 keeps the leading zero as text; C11 is a literal
# number 974123658 but still carries the text style index.
$ws.Cells.Item(11,1).NumberFormat = "@"
$ws.Cells.Item(11,1).Value = "0987456231"
$ws.Cells.Item(11,2).Value = "Diego Montalvo"
$ws.Cells.Item(11,3).Value = 974123658
$ws.Cells.Item(11,3).NumberFormat = "@"
$ws.Cells.Item(11,4).Value = 27
$ws.Cells.Item(11,5).Value = "ADPR321000"

$ws.Range("B11").Select()
